$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 142.9073533333333
$ws.Range("H2").Value = 428.72206
$ws.Range("I2").Value = 0.5576664151504187
$ws.Range("J2").Value = 0.5576664151504188
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.324764666666667
$ws.Range("N2").Value = 6.974294
$ws.Range("O2").Value = 0.04473923998638302
$ws.Range("P2").Value = 0.04473923998638301
$ws.Range("Q2").Value = 332.2259656361823
$ws.Range("R2").Value = 2990.03369072564
$ws.Range("S2").Value = 0.02494957157976048
$ws.Range("T2").Value = 0.02494957157976049

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 142.9073533333333
$ws.Range("H3").Value = 428.72206
$ws.Range("I3").Value = 0.5576664151504187
$ws.Range("J3").Value = 0.5576664151504188
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 19.27491966666667
$ws.Range("N3").Value = 57.824759
$ws.Range("O3").Value = 0.3709387315842666
$ws.Range("P3").Value = 0.3709387315842665
$ws.Range("Q3").Value = 2754.527755275949
$ws.Range("R3").Value = 24790.74979748354
$ws.Range("S3").Value = 0.2068600726830414
$ws.Range("T3").Value = 0.2068600726830414

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 142.9073533333333
$ws.Range("H4").Value = 428.72206
$ws.Range("I4").Value = 0.5576664151504187
$ws.Range("J4").Value = 0.5576664151504188
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 30.36285833333334
$ws.Range("N4").Value = 91.088575
$ws.Range("O4").Value = 0.5843220284293504
$ws.Range("P4").Value = 0.5843220284293504
$ws.Range("Q4").Value = 4339.075724051611
$ws.Range("R4").Value = 39051.6815164645
$ws.Range("S4").Value = 0.3258567708876169
$ws.Range("T4").Value = 0.3258567708876169

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 63.967809
$ws.Range("H5").Value = 191.903427
$ws.Range("I5").Value = 0.2496211559306514
$ws.Range("J5").Value = 0.2496211559306514
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.324764666666667
$ws.Range("N5").Value = 6.974294
$ws.Range("O5").Value = 0.04473923998638302
$ws.Range("P5").Value = 0.04473923998638301
$ws.Range("Q5").Value = 148.710102167282
$ws.Range("R5").Value = 1338.390919505538
$ws.Range("S5").Value = 0.01116786080085975
$ws.Range("T5").Value = 0.01116786080085975

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 63.967809
$ws.Range("H6").Value = 191.903427
$ws.Range("I6").Value = 0.2496211559306514
$ws.Range("J6").Value = 0.2496211559306514
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 19.27491966666667
$ws.Range("N6").Value = 57.824759
$ws.Range("O6").Value = 0.3709387315842666
$ws.Range("P6").Value = 0.3709387315842665
$ws.Range("Q6").Value = 1232.974379727677
$ws.Range("R6").Value = 11096.76941754909
$ws.Range("S6").Value = 0.09259415495751425
$ws.Range("T6").Value = 0.09259415495751425

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 63.967809
$ws.Range("H7").Value = 191.903427
$ws.Range("I7").Value = 0.2496211559306514
$ws.Range("J7").Value = 0.2496211559306514
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 30.36285833333334
$ws.Range("N7").Value = 91.088575
$ws.Range("O7").Value = 0.5843220284293504
$ws.Range("P7").Value = 0.5843220284293504
$ws.Range("Q7").Value = 1942.245522560725
$ws.Range("R7").Value = 17480.20970304653
$ws.Range("S7").Value = 0.1458591401722774
$ws.Range("T7").Value = 0.1458591401722774

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 49.38440333333333
$ws.Range("H8").Value = 148.15321
$ws.Range("I8").Value = 0.1927124289189298
$ws.Range("J8").Value = 0.1927124289189298
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.324764666666667
$ws.Range("N8").Value = 6.974294
$ws.Range("O8").Value = 0.04473923998638302
$ws.Range("P8").Value = 0.04473923998638301
$ws.Range("Q8").Value = 114.8071159537489
$ws.Range("R8").Value = 1033.26404358374
$ws.Range("S8").Value = 0.008621807605762779
$ws.Range("T8").Value = 0.008621807605762779

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 49.38440333333333
$ws.Range("H9").Value = 148.15321
$ws.Range("I9").Value = 0.1927124289189298
$ws.Range("J9").Value = 0.1927124289189298
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 19.27491966666667
$ws.Range("N9").Value = 57.824759
$ws.Range("O9").Value = 0.3709387315842666
$ws.Range("P9").Value = 0.3709387315842665
$ws.Range("Q9").Value = 951.8804070362654
$ws.Range("R9").Value = 8566.92366332639
$ws.Range("S9").Value = 0.07148450394371096
$ws.Range("T9").Value = 0.07148450394371096

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 49.38440333333333
$ws.Range("H10").Value = 148.15321
$ws.Range("I10").Value = 0.1927124289189298
$ws.Range("J10").Value = 0.1927124289189298
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 30.36285833333334
$ws.Range("N10").Value = 91.088575
$ws.Range("O10").Value = 0.5843220284293504
$ws.Range("P10").Value = 0.5843220284293504
$ws.Range("Q10").Value = 1499.451642286195
$ws.Range("R10").Value = 13495.06478057575
$ws.Range("S10").Value = 0.1126061173694561
$ws.Range("T10").Value = 0.1126061173694561

